$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(29)
$rng = $p.Range
Write-Output ("before text=[" + $rng.Text + "]")
$xml = @'
<w:p w14:paraId="26663E9A" w14:textId="4920AC89" w:rsidR="000F1AF2" w:rsidRPr="003A3DEB" w:rsidRDefault="000F1AF2" w:rsidP="003A3DEB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r w:rsidRPr="003A3DEB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Liu, X, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003A3DEB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Takecuchi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="003A3DEB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">, Y., </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003A3DEB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Iwami</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="003A3DEB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>, S., A Mathematical Investigation of Vaccination Strategies to Prevent Measles Epidemics. Journal of Theoretical Biology 253 (2008) 1&#8211;11.</w:t></w:r></w:p>
'@
$rng.InsertXML($xml)
Write-Output ("Para count: " + $d.Paragraphs.Count)
